$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the "Periodo Mora" / "Valor Mora" table (rows 16-19) into
# ascending period order (2202, 2203, 2204, 2205), keeping each
# period's corresponding "Valor Mora" amount attached to it.
$periods = @("2202", "2203", "2204", "2205")
$valores = @(40000, 40000, 40000, 6667)

for ($i = 0; $i -lt 4; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
